$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.348.13"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "2.622.12"
$ws.Range("E3").Value = "  -2.08%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'595.33"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").Value = "'167.28"
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -2.29%  "
$ws.Range("D9").Value = "2.621.97"
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("D12").Value = "'0.363"
$ws.Range("E12").Value = "  +1.60%  "
$ws.Range("D13").Value = "'5.24"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").Value = "'27.70"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "3.098.86"
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D17").Value = "67.040.92"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").Value = "2.620.63"
$ws.Range("E18").Value = "  -2.90%  "
$ws.Range("D19").Value = "'12.08"
$ws.Range("E19").Value = "  +3.05%  "
$ws.Range("D20").Value = "'7.96"
$ws.Range("E20").Value = "  +4.61%  "
$ws.Range("D21").Value = "'357.49"
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("D22").Value = "'4.32"
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("D23").Value = "'4.66"
$ws.Range("E23").Value = "  -3.21%  "
$ws.Range("D25").Value = "'1.93"
$ws.Range("E25").Value = "  -5.18%  "
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").Value = "'69.72"
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("D28").Value = "2.740.72"
$ws.Range("E28").Value = "  -2.33%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").Value = "'0.0000100"
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("D31").Value = "'544.27"
$ws.Range("E31").Value = "  -1.91%  "
$ws.Range("D32").Value = "'7.91"
$ws.Range("E32").Value = "  -1.12%  "
$ws.Range("E33").Value = "  -2.82%  "
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("D35").Value = "'0.136"
$ws.Range("E35").Value = "  +5.05%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  -3.83%  "
$ws.Range("D38").Value = "'157.20"
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("D39").Value = "'18.99"
$ws.Range("E39").Value = "  -2.79%  "
$ws.Range("E40").Value = "  -2.02%  "
$ws.Range("D41").Value = "'18.18"
$ws.Range("E41").Value = "  +1.36%  "
$ws.Range("D42").Value = "'1.81"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").Value = "'5.21"
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("E45").Value = "  -4.25%  "
$ws.Range("D46").Value = "0.0₆0297"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").Value = "'151.97"
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("D48").Value = "'0.579"
$ws.Range("E48").Value = "  -2.23%  "
$ws.Range("E49").Value = "  -1.84%  "
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").Value = "'0.0770"
$ws.Range("E51").Value = "  -0.94%  "

# Reset style back to Normal for cells that needed a quote-prefix to stay text,
# so no stray number-format/style index is left on the cell.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Style = "Normal"
